$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.24"
$ws.Range("E2").Value = "'1.81%"
$ws.Range("D3").Value = "'42.25"
$ws.Range("E3").Value = "'4.31%"
$ws.Range("D4").Value = "'5.010"
$ws.Range("E4").Value = "'-0.10%"
$ws.Range("D5").Value = "'0.07521"
$ws.Range("E5").Value = "'3.05%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.368"
$ws.Range("E6").Value = "'2.02%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.596"
$ws.Range("E7").Value = "'2.07%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9309"
$ws.Range("E8").Value = "'0.19%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.401"
$ws.Range("E9").Value = "'2.92%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1190"
$ws.Range("E10").Value = "'2.46%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1835"
$ws.Range("E11").Value = "'4.73%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08906"
$ws.Range("E12").Value = "'2.11%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04134"
$ws.Range("E13").Value = "'-5.37%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1050"
$ws.Range("E14").Value = "'-0.31%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001292"
$ws.Range("E15").Value = "'2.22%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005800"
$ws.Range("E16").Value = "'-3.67%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.336"
$ws.Range("E17").Value = "'-0.20%"
$ws.Range("D18").Value = "'0.3333"
$ws.Range("D19").Value = "'8.317"
$ws.Range("E19").Value = "'5.45%"
$ws.Range("E20").Value = "'-2.81%"
$ws.Range("D22").Value = "'0.04083"
$ws.Range("E22").Value = "'4.14%"
$ws.Range("E23").Value = "'0.39%"
$ws.Range("D24").Value = "'0.003892"
$ws.Range("E24").Value = "'5.67%"
$ws.Range("D25").Value = "'0.0001300"
$ws.Range("E25").Value = "'8.35%"
$ws.Range("D38").Value = "'0.02395"
$ws.Range("E38").Value = "'3.52%"
$ws.Range("D39").Value = "'0.05223"
$ws.Range("E39").Value = "'2.93%"
$ws.Range("D40").Value = "'0.006749"
$ws.Range("E40").Value = "'16.74%"
$ws.Range("D41").Value = "'0.007770"
$ws.Range("E41").Value = "'-0.90%"
$ws.Range("E42").Value = "'3.03%"
$ws.Range("D43").Value = "'0.007407"
$ws.Range("E43").Value = "'0.35%"
$ws.Range("D44").Value = "'0.007123"
$ws.Range("E44").Value = "'-1.44%"
$ws.Range("D45").Value = "'0.2986"
$ws.Range("E45").Value = "'2.38%"
$ws.Range("D46").Value = "'0.00006568"
$ws.Range("E46").Value = "'5.86%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.03%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.004203"
$ws.Range("E48").Value = "'0.03%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.04519"
$ws.Range("E49").Value = "'-6.48%"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E51").Value = "'-0.03%"
